# Add the missing SLEEP_CONTROL register (with its bit fields) to the
# WOLFE_APB_SOC_CTRL reference workbook.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # IPREGLIST_rel1.0.0
$ws2 = $wb.Worksheets.Item(2)   # IPREGMAP_rel1.0.0

# ---------------------------------------------------------------------
# Sheet 1 (IPREGLIST_rel1.0.0): insert a new register row (SLEEP_CONTROL)
# above the existing RTC_FIRST_REG / RTC_LAST_REG rows, so the existing
# two data rows shift down from 2,3 to 3,4.
# ---------------------------------------------------------------------

$ws1.Rows.Item(2).Insert()

$row = 2
$ws1.Cells.Item($row, 1).Value2 = "SLEEP_CONTROL"
$ws1.Cells.Item($row, 3).Value2 = "0x104"
$ws1.Cells.Item($row, 4).Value2 = 32
$ws1.Cells.Item($row, 5).Value2 = "Config"
$ws1.Cells.Item($row, 6).Value2 = "R/W"
$ws1.Cells.Item($row, 7).Value2 = "PERIPH"
$ws1.Cells.Item($row, 8).Value2 = "R/W"
$ws1.Cells.Item($row, 9).Value2 = "0x0"
$ws1.Cells.Item($row, 10).Value2 = "Deep sleep control register"

# Match the look of the surrounding data rows (center aligned, wrapped
# text) for the newly inserted row.
$dataRange = $ws1.Range("A2:J2")
$dataRange.HorizontalAlignment = -4108
$dataRange.WrapText = $true
$ws1.Range("B2:K2").HorizontalAlignment = -4108
$ws1.Range("B2,K2").WrapText = $false
$ws1.Rows.Item(2).RowHeight = 14.25
$ws1.Rows.Item(3).RowHeight = 13.8
$ws1.Rows.Item(4).RowHeight = 13.8

# ---------------------------------------------------------------------
# Sheet 2 (IPREGMAP_rel1.0.0): add the bit-field rows describing the new
# SLEEP_CONTROL register.
# ---------------------------------------------------------------------

function Set-BitField {
    param($r, $name, $register, $bitPos, $bitSize, $desc, $wrap)

    $ws2.Cells.Item($r, 1).Value2 = $name
    $ws2.Cells.Item($r, 2).Value2 = $register
    $ws2.Cells.Item($r, 3).Value2 = $bitPos
    $ws2.Cells.Item($r, 4).Value2 = $bitSize
    $ws2.Cells.Item($r, 5).Value2 = "R/W"
    $ws2.Cells.Item($r, 6).Value2 = "R/W"
    $ws2.Cells.Item($r, 7).Value2 = "0x0"
    $ws2.Cells.Item($r, 8).Value2 = $desc

    $ws2.Range($ws2.Cells.Item($r, 1), $ws2.Cells.Item($r, 1)).HorizontalAlignment = -4108
    $ws2.Range($ws2.Cells.Item($r, 2), $ws2.Cells.Item($r, 7)).HorizontalAlignment = -4108
    $ws2.Range($ws2.Cells.Item($r, 2), $ws2.Cells.Item($r, 2)).WrapText = $true
    $ws2.Range($ws2.Cells.Item($r, 5), $ws2.Cells.Item($r, 7)).WrapText = $true
    $ws2.Cells.Item($r, 8).HorizontalAlignment = -4131
    $ws2.Cells.Item($r, 8).WrapText = $wrap
}

$extWakeupType = "External wakeup type. This tells the way the external GPIO can wakeup the chip while it is in deep sleep (raising edge, falling edge, etc). Possible values:`n- 0: Rising edge.`n- 1: Falling edge.`n- 2: Level high.`n- 3: Level low."

Set-BitField 2 "FLL_RET"        "SLEEP_CONTROL" 0  2 "FLL retention configuration." $true
Set-BitField 3 "MEM_RET_0"      "SLEEP_CONTROL" 2  1 "Memory retention configuration." $true
Set-BitField 4 "EXTWAKEUP_SEL"  "SLEEP_CONTROL" 6  5 "External wakeup selection. This gives the GPIO numer which can wakeup the chip when it is in deep sleep mode." $false
Set-BitField 5 "EXTWAKEUP_TYPE" "SLEEP_CONTROL" 11 2 $extWakeupType $true
Set-BitField 6 "EXTWAKEUP_EN"   "SLEEP_CONTROL" 13 1 "External wakeup enable." $false

$ws2.Rows.Item(2).RowHeight = 14.25
$ws2.Rows.Item(3).RowHeight = 14.25
$ws2.Rows.Item(5).RowHeight = 66.1

# ---------------------------------------------------------------------
# Final view state: sheet 2 (IPREGMAP) becomes the active tab, with a
# slightly reduced zoom level on both sheets, and updated selections.
# ---------------------------------------------------------------------

$ws1.Activate()
$ws1.Range("J3").Select()
$excel.ActiveWindow.Zoom = 55

$ws2.Activate()
$ws2.Range("H7").Select()
$excel.ActiveWindow.Zoom = 55
